# Apply updates to column F (dSF) for a set of rows, per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    7  = 1
    12 = 0
    17 = 2
    25 = -4
    26 = 0
    33 = -2
    34 = -2
    36 = 4
    41 = 2
    46 = 2
    51 = 2
    52 = 0
    55 = 2
    56 = -1
    62 = 0
    63 = -10
    64 = 10
    66 = -6
    70 = -6
    74 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
